# Applies a weekly refresh of the "Fruta / hortaliza" data: the values in
# columns D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado) and P (Precio $/Kg) are re-shuffled across
# rows 2-16 (row 8 is left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Maps destination row -> source row (values to copy from the row's
# current/original contents before any cell is overwritten).
$rowMap = @{
    2  = 15
    3  = 12
    4  = 3
    5  = 9
    6  = 14
    7  = 13
    8  = 8
    9  = 10
    10 = 11
    11 = 2
    12 = 5
    13 = 6
    14 = 16
    15 = 4
    16 = 7
}

$cols = @("D", "J", "K", "L", "M", "P")

# Snapshot the original values for the columns we are about to shuffle so
# that overwriting one row doesn't corrupt the source data for another.
$original = @{}
foreach ($row in $rowMap.Keys) {
    $rowValues = @{}
    foreach ($col in $cols) {
        $rowValues[$col] = $ws.Range("$col$row").Value()
    }
    $original[$row] = $rowValues
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $original[$srcRow][$col]
    }
}
